# Apply updated odds values to the "Jogos da Semana" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 changes
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 11
$ws.Range("Q3").Value = 1.95
$ws.Range("R3").Value = 1.9
$ws.Range("W3").Value = 6.5
$ws.Range("Z3").Value = 9.5
$ws.Range("AD3").Value = 8
$ws.Range("AG3").Value = 451
$ws.Range("AI3").Value = 34
$ws.Range("AJ3").Value = 21
$ws.Range("BC3").Value = 126

# Row 5 changes
$ws.Range("G5").Value = 2.63
$ws.Range("I5").Value = 2.9
$ws.Range("J5").Value = 3.4
$ws.Range("L5").Value = 3.75
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("W5").Value = 6.5
$ws.Range("X5").Value = 11
$ws.Range("AI5").Value = 13
$ws.Range("AO5").Value = 15
$ws.Range("AY5").Value = 34

# Row 6 changes
$ws.Range("I6").Value = 2.6
$ws.Range("N6").Value = 7.3
$ws.Range("W6").Value = 8.5
$ws.Range("Y6").Value = 10
$ws.Range("AE6").Value = 13.5
$ws.Range("AH6").Value = 7.4
$ws.Range("AI6").Value = 12.5
$ws.Range("AK6").Value = 30
$ws.Range("AL6").Value = 24
